$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift the key/value (columns B and C) of rows 99..111 down by one row,
# working from the bottom up so we never clobber a source cell before
# it has been read. This makes room for a brand-new entry at row 99:
#   XML2SOURCE_FILE.REQUEST_FACTORY_RETURN.DESCRIPTION / "{0} クラスを返却します。"
for ($r = 111; $r -ge 99; $r--) {
    $srcB = $ws.Cells.Item($r, 2).Value2
    $srcC = $ws.Cells.Item($r, 3).Value2

    if ($srcB -eq $null) { $srcB = "" }
    if ($srcC -eq $null) { $srcC = "" }

    $ws.Cells.Item($r + 1, 2).Value = $srcB
    $ws.Cells.Item($r + 1, 3).Value = $srcC
}

# The table grew by one row, so the two trailing template rows shuffle
# down: old row 113 (the final/closing row) moves to 114, and old row 112
# (the blank "spacer" row) moves to 113 -- copy formatting bottom-up so
# we don't clobber a source row's look before it has been duplicated.
$ws.Range("A113:G113").Copy()
$ws.Range("A114:G114").PasteSpecial(-4122)
$ws.Range("A112:G112").Copy()
$ws.Range("A113:G113").PasteSpecial(-4122)

# Row 112 itself becomes a real data row, so give it the same cell
# formatting (borders etc.) as the data rows above it.
$ws.Range("A111:G111").Copy()
$ws.Range("A112:G112").PasteSpecial(-4122)

# Clear any leftover values those rows might have had (they are blank
# template rows), but keep their (freshly pasted) formatting.
$ws.Range("A113:G113").ClearContents()
$ws.Range("A114:G114").ClearContents()

# The running "No." counter in column A is a fill-down formula
# (=A<previous row>+1); extend it onto the newly-activated row 112.
$ws.Range("A112").Formula = "=A111+1"

# Fill in the new row 99 with the new resource-bundle key and message.
$ws.Cells.Item(99, 2).Value = "XML2SOURCE_FILE.REQUEST_FACTORY_RETURN.DESCRIPTION"
$ws.Cells.Item(99, 3).Value = "{0} クラスを返却します。"

# Keep the selection roughly where the author left it.
$ws.Range("C100").Select() | Out-Null
